$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClickThroughRateSheet")

# Touching these previously-untouched cells with the default "Normal" style
# materializes them as empty <c/> placeholders in the saved sheet XML,
# expanding the sheet's recorded used range to A1:D13 without altering any
# existing value/style.
$ws.Range("B1:D1").Style = "Normal"
$ws.Range("A2:A13").Style = "Normal"

# Refresh the sample data (Clicks / Impressions) that feeds the CTR formulas.
$ws.Range("B3").Value = 509
$ws.Range("C3").Value = 346380

$ws.Range("B4").Value = 614
$ws.Range("C4").Value = 293935

$ws.Range("B5").Value = 70
$ws.Range("C5").Value = 899517

$ws.Range("B6").Value = 617
$ws.Range("C6").Value = 781247

$ws.Range("B7").Value = 975
$ws.Range("C7").Value = 800261

$ws.Range("B8").Value = 99
$ws.Range("C8").Value = 333038

$ws.Range("B9").Value = 810
$ws.Range("C9").Value = 843502

$ws.Range("B10").Value = 976
$ws.Range("C10").Value = 132730

$ws.Range("B11").Value = 214
$ws.Range("C11").Value = 57169

$ws.Range("B12").Value = 961
$ws.Range("C12").Value = 717649

$ws.Range("B13").Value = 880
$ws.Range("C13").Value = 32781
